$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2147.2
$ws.Range("I86").Value = 2166.3333
$ws.Range("J86").Value = 2118.5
$ws.Range("K86").Value = 2166.3333
$ws.Range("L86").Value = 2118.5
$ws.Range("M86").Value = -1043.3333
$ws.Range("N86").Value = -4364.5
$ws.Range("H89").Value = 2147.2
$ws.Range("I89").Value = 2166.3333
$ws.Range("J89").Value = 2118.5
$ws.Range("K89").Value = 10831.6665
$ws.Range("L89").Value = 10592.5
$ws.Range("M89").Value = -5215.666499999999
$ws.Range("N89").Value = -21824.5
$ws.Range("H116").Value = 8651.833000000001
$ws.Range("I116").Value = 7635
$ws.Range("J116").Value = 9668.666999999999
$ws.Range("K116").Value = 7635
$ws.Range("L116").Value = 9668.666999999999
$ws.Range("M116").Value = -4193
$ws.Range("N116").Value = -16552.667
$ws.Range("H125").Value = 5498.1665
$ws.Range("J125").Value = 5397.8
$ws.Range("L125").Value = 48580.2
$ws.Range("N125").Value = -53500.2
$ws.Range("H135").Value = 1960.3914
$ws.Range("I135").Value = 752.2941
$ws.Range("J135").Value = 5383.3335
$ws.Range("K135").Value = 6770.6469
$ws.Range("L135").Value = 48450.0015
$ws.Range("M135").Value = -4235.6469
$ws.Range("N135").Value = -53520.0015
$ws.Range("H138").Value = 2695.5122
$ws.Range("I138").Value = 1911.2
$ws.Range("J138").Value = 3148
$ws.Range("K138").Value = 5733.6
$ws.Range("L138").Value = 9444
$ws.Range("M138").Value = -593.6000000000004
$ws.Range("N138").Value = -19724

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5348574
$ws.Range("I2").Value = 6493868.5
$ws.Range("J2").Value = 3866.6667
$ws.Range("K2").Value = 6493868.5
$ws.Range("L2").Value = 3866.6667
$ws.Range("M2").Value = -6493755.5
$ws.Range("N2").Value = -4092.6667
$ws.Range("H6").Value = 19800
$ws.Range("J6").Value = 19750
$ws.Range("L6").Value = 19750
$ws.Range("N6").Value = -20096
$ws.Range("H32").Value = 27103.633
$ws.Range("I32").Value = 27132.152
$ws.Range("K32").Value = 27132.152
$ws.Range("M32").Value = -26845.152
$ws.Range("H43").Value = 31475
$ws.Range("I43").Value = 19000
$ws.Range("K43").Value = 19000
$ws.Range("M43").Value = -18687
$ws.Range("H45").Value = 1855
$ws.Range("I45").Value = 1746.6666
$ws.Range("K45").Value = 1746.6666
$ws.Range("M45").Value = -1369.6666
$ws.Range("H74").Value = 2236.8667
$ws.Range("I74").Value = 1173.4736
$ws.Range("K74").Value = 1173.4736
$ws.Range("M74").Value = -299.4736
$ws.Range("H77").Value = 2236.8667
$ws.Range("I77").Value = 1173.4736
$ws.Range("K77").Value = 5867.368
$ws.Range("M77").Value = -1499.368
$ws.Range("H110").Value = 13890800
$ws.Range("I110").Value = 19231876
$ws.Range("K110").Value = 19231876
$ws.Range("M110").Value = -19229831
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 5348574
$ws.Range("I116").Value = 6493868.5
$ws.Range("J116").Value = 3866.6667
$ws.Range("K116").Value = 6493868.5
$ws.Range("L116").Value = 3866.6667
$ws.Range("M116").Value = -6491574.5
$ws.Range("N116").Value = -8454.6667
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 4462.58
$ws.Range("I132").Value = 3164.4033
$ws.Range("K132").Value = 9493.2099
$ws.Range("M132").Value = -6963.2099

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5348574
$ws.Range("I3").Value = 6493868.5
$ws.Range("J3").Value = 3866.6667
$ws.Range("K3").Value = 6493868.5
$ws.Range("L3").Value = 3866.6667
$ws.Range("M3").Value = -6493754.5
$ws.Range("N3").Value = -4094.6667
$ws.Range("H107").Value = 1251.2413
$ws.Range("I107").Value = 1221.1666
$ws.Range("K107").Value = 1221.1666
$ws.Range("M107").Value = 698.8334
$ws.Range("H134").Value = 5766.878
$ws.Range("I134").Value = 4134.3794
$ws.Range("K134").Value = 12403.1382
$ws.Range("M134").Value = -9868.138199999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 35.142857
$ws.Range("I7").Value = 33.666668
$ws.Range("J7").Value = 36.25
$ws.Range("K7").Value = 33.666668
$ws.Range("L7").Value = 36.25
$ws.Range("M7").Value = 79.333332
$ws.Range("N7").Value = -262.25
$ws.Range("H31").Value = 26319446
$ws.Range("J31").Value = 4057.375
$ws.Range("L31").Value = 4057.375
$ws.Range("N31").Value = -4647.375
$ws.Range("H34").Value = 26319446
$ws.Range("J34").Value = 4057.375
$ws.Range("L34").Value = 4057.375
$ws.Range("N34").Value = -4461.375
$ws.Range("H35").Value = 4875
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20588
$ws.Range("H81").Value = 59999
$ws.Range("J81").Value = 59999
$ws.Range("L81").Value = 59999
$ws.Range("N81").Value = -61995
$ws.Range("H84").Value = 59999
$ws.Range("J84").Value = 59999
$ws.Range("L84").Value = 179997
$ws.Range("N84").Value = -189981
$ws.Range("H96").Value = 25543.8
$ws.Range("J96").Value = 25543.8
$ws.Range("L96").Value = 25543.8
$ws.Range("N96").Value = -31035.8
$ws.Range("H106").Value = 48556.668
$ws.Range("J106").Value = 48556.668
$ws.Range("L106").Value = 48556.668
$ws.Range("N106").Value = -51080.668
$ws.Range("H132").Value = 38761.812
$ws.Range("I132").Value = 3446.7856
$ws.Range("K132").Value = 10340.3568
$ws.Range("M132").Value = -7810.356800000001
$ws.Range("H134").Value = 4746.421
$ws.Range("I134").Value = 2888.0715
$ws.Range("K134").Value = 8664.2145
$ws.Range("M134").Value = -6129.2145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3453.5454
$ws.Range("J62").Value = 3250
$ws.Range("L62").Value = 9750
$ws.Range("N62").Value = -11122
$ws.Range("H65").Value = 3453.5454
$ws.Range("J65").Value = 3250
$ws.Range("L65").Value = 29250
$ws.Range("N65").Value = -36114
$ws.Range("H107").Value = 1433.75
$ws.Range("I107").Value = 410.88235
$ws.Range("K107").Value = 1232.64705
$ws.Range("M107").Value = 687.35295

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 59237.5
$ws.Range("J32").Value = 69000
$ws.Range("L32").Value = 69000
$ws.Range("N32").Value = -69592
$ws.Range("H102").Value = 1170
$ws.Range("I102").Value = 1170
$ws.Range("K102").Value = 1170
$ws.Range("M102").Value = 452
$ws.Range("H113").Value = 3801338.2
$ws.Range("I113").Value = 1212
$ws.Range("K113").Value = 1212
$ws.Range("M113").Value = 958
$ws.Range("H126").Value = 3477.238
$ws.Range("I126").Value = 2481.4546
$ws.Range("J126").Value = 4572.6
$ws.Range("K126").Value = 7444.3638
$ws.Range("L126").Value = 13717.8
$ws.Range("M126").Value = -4974.3638
$ws.Range("N126").Value = -18657.8
$ws.Range("H129").Value = 21990
$ws.Range("J129").Value = 21990
$ws.Range("L129").Value = 21990
$ws.Range("N129").Value = -31990

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("H55").Value = 142.16667
$ws.Range("J55").Value = 172.5
$ws.Range("L55").Value = 172.5
$ws.Range("N55").Value = -518.5
$ws.Range("H104").Value = 19998.334
$ws.Range("J104").Value = 19998.334
$ws.Range("L104").Value = 19998.334
$ws.Range("N104").Value = -26986.334
$ws.Range("H122").Value = 4227.909
$ws.Range("I122").Value = 3610
$ws.Range("K122").Value = 10830
$ws.Range("M122").Value = -8380
$ws.Range("H126").Value = 3000

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 36044.8
$ws.Range("J95").Value = 36044.8
$ws.Range("L95").Value = 36044.8
$ws.Range("N95").Value = -41536.8
$ws.Range("H103").Value = 53398.6
$ws.Range("J103").Value = 53398.6
$ws.Range("L103").Value = 53398.6
$ws.Range("N103").Value = -55742.6
$ws.Range("H126").Value = 8778.571
$ws.Range("I126").Value = 10611.25
$ws.Range("J126").Value = 6335
$ws.Range("K126").Value = 31833.75
$ws.Range("L126").Value = 19005
$ws.Range("M126").Value = -29363.75
$ws.Range("N126").Value = -23945
$ws.Range("H132").Value = 5068.8965
$ws.Range("I132").Value = 3439
$ws.Range("J132").Value = 9347.375
$ws.Range("K132").Value = 10317
$ws.Range("L132").Value = 28042.125
$ws.Range("M132").Value = -7787
$ws.Range("N132").Value = -33102.125
$ws.Range("H136").Value = 4294.909
$ws.Range("J136").Value = 5129.1875
$ws.Range("L136").Value = 15387.5625
$ws.Range("N136").Value = -20487.5625
